# Append: 2025-10-17 06:25 JST
# Re-scrape refreshed the "ランサーズ" listing sheet: new snapshot timestamp,
# only 4 surviving/new listings (old rows 5, 6 and 17 survive, shifted up to
# rows 2-4 are replaced by one brand-new listing, row 5 keeps the old row17),
# column B/D narrowed, and stale rows + their hyperlinks are dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2025-10-17 06:24:59"

# --- Column widths (ColumnWidth is quantized to 1/7-character steps by the
# engine's pixel round-trip, so nudge by 1/7 to land exactly on the target
# integer width as stored in the OOXML <col width="..."/>). ---
$ws.Columns.Item(2).ColumnWidth = 37 + 1/7   # B: 52 -> 38
$ws.Columns.Item(4).ColumnWidth = 27 + 1/7   # D: 32 -> 28

# --- Drop every stale hyperlink first (rebuilt below for the surviving
# rows only); Hyperlinks.Delete() clears the whole sheet collection. ---
$ws.Hyperlinks.Delete()

# --- Row 2: RoboTANGO (was row 5) ---
$ws.Cells.Item(2, 1).Value = $timestamp
$ws.Cells.Item(2, 2).Value = "【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5405023"
$ws.Cells.Item(2, 7).Value = 178
$ws.Cells.Item(2, 8).Value = "★bot ◆ツール"

# --- Row 3: 在庫管理・出品補助ツール Zoom面談依頼 (was row 6) ---
$ws.Cells.Item(3, 1).Value = $timestamp
$ws.Cells.Item(3, 2).Value = "【相談希望】在庫管理・出品補助ツールの開発に関するZoom面談依頼"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5398112"
$ws.Cells.Item(3, 7).Value = 158
$ws.Cells.Item(3, 8).Value = "◆ツール,開発 ◇管理"

# --- Row 4: new listing, 医薬品マッチング (no skill-tag column) ---
$ws.Cells.Item(4, 1).Value = $timestamp
$ws.Cells.Item(4, 2).Value = "【医薬品マッチング】高額医薬品の譲渡支援システム構築"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5415061"
$ws.Cells.Item(4, 7).Value = 33
$ws.Cells.Item(4, 8).ClearContents()

# --- Row 5: VBA 1問1答問題集 (was row 17) ---
$ws.Cells.Item(5, 1).Value = $timestamp
$ws.Cells.Item(5, 2).Value = "初回 【急募・即決します】VBAで1問1答問題集の作成"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5414812"
$ws.Cells.Item(5, 7).Value = 10
$ws.Cells.Item(5, 8).ClearContents()

# --- Drop the now-stale rows 6:17 entirely so the used range / dimension
# shrinks back down to A1:H5. ---
$ws.Range("A6:H17").EntireRow.Delete()

# --- Re-create hyperlinks only for the 4 surviving URL cells. ---
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.lancers.jp/work/detail/5405023")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.lancers.jp/work/detail/5398112")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.lancers.jp/work/detail/5415061")
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5414812")
